$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Picture run (paragraph with "Het eerste punt leek ons...") gains a
#    <w:lang w:eastAsia="nl-NL"/> entry in its rPr (alongside the existing
#    <w:noProof/>). We round-trip that single paragraph's OOXML, patch the
#    one run whose rPr is bare "<w:noProof/>", strip the synthetic
#    w14:paraId/w14:textId that the getter stamps on, and feed it back in.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "*Het eerste punt leek ons geen probleem*") {
        $targetPara = $cand
        break
    }
}

$picXml = $targetPara.Range.WordOpenXML
$picXml = $picXml.Replace('<w:rPr><w:noProof/></w:rPr><w:drawing>', '<w:rPr><w:noProof/><w:lang w:eastAsia="nl-NL"/></w:rPr><w:drawing>')
$picXml = $picXml -replace ' w14:paraId="[0-9A-Fa-f]+" w14:textId="[0-9A-Fa-f]+"', ''
$targetPara.Range.InsertXML($picXml)

# ---------------------------------------------------------------------------
# 2) Locate the paragraph that currently holds the "_GoBack" bookmark plus
#    the "De vertraging key ..." sentence. Extend it in place with the new
#    "if-statement" explanation text (the original runs + proofErr markers
#    are hand authored here so they are not lost/merged by a WordOpenXML
#    round trip), dropping the bookmark (it moves further down below).
# ---------------------------------------------------------------------------
$delayPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "De vertraging key word alleen meegegeven*") {
        $delayPara = $cand
        break
    }
}

$delayInner = '<w:r><w:t xml:space="preserve">De vertraging </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>key</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> word alleen meegegeven als er vertraging is echter was deze ingesteld als een waarde die niet nul kon zijn. </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">Dit probleem is opgelost door een </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>if</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve">-statement toe te voegen die controleer of de </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>key</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> bestaat, als dit het geval is word deze waarde meegegeven. Als de </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>key</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> niet bestaat word er een 0 waarde aangemaakt en meegegeven.</w:t></w:r>'

$delayXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p>' + $delayInner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$delayPara.Range.InsertXML($delayXml)

# ---------------------------------------------------------------------------
# 3) Append four new paragraphs after it:
#      - the "niet statisch station" / zoekveld paragraph
#      - a "Github" Heading 2
#      - the "Voor we begonnen met git" paragraph
#      - the closing "Dit was dan ook voor ons het keerpunt..." paragraph,
#        which is where the _GoBack bookmark now lives.
# ---------------------------------------------------------------------------
$delayPara.Range.InsertParagraphAfter()
$searchPara = $d.Paragraphs($delayPara.Index + 1)

$searchInner = '<w:r><w:t xml:space="preserve">Nadat ook dit probleem was opgelost zijn we aan de slag gegaan met een niet statisch station. Voorheen kon er nog niet op een station gezocht worden en zaten we dus vast aan de gegevens van een vast station. Hiervoor is een zoekveld toegevoegd aan de app, waarbij de </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r><w:t>tableview</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> zich steeds vernieuwt als er een verandering plaats vind van de tekst in dit veld.</w:t></w:r>'

$searchXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p>' + $searchInner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$searchPara.Range.InsertXML($searchXml)

$searchPara.Range.InsertParagraphAfter()
$githubPara = $d.Paragraphs($searchPara.Index + 1)

$githubInner = '<w:pPr><w:pStyle w:val="Kop2"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Github</w:t></w:r><w:proofErr w:type="spellEnd"/>'

$githubXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p>' + $githubInner + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$githubPara.Range.InsertXML($githubXml)

$githubPara.Range.InsertParagraphAfter()
$gitParaA = $d.Paragraphs($githubPara.Index + 1)
$gitParaA.Range.InsertAfter("Voor we begonnen met het gebruik van git zaten we steeds vast op een laptop. Waarbij we dan ook met twee man naar een scherm moesten staren. Dit was alles behalve effici" + [char]0x00EB + "nt. Als gevolg hiervan was dan ook dat we veel tijd kwijt waren met vrij simpele dingen. Maar ook het versie beheer werd niet gedaan. Dit heeft er dan ook voor gezorgd dat we een keer een stuk code hebben moeten herschrijven omdat deze het niet meer deed na het toevoegen van een ander stuk code.")

$gitParaA.Range.InsertParagraphAfter()
$gitParaB = $d.Paragraphs($gitParaA.Index + 1)
$gitParaB.Range.InsertAfter("Dit was dan ook voor ons het keerpunt waarbij we begonnen zijn met git. ")

# Move the _GoBack bookmark to the end of this final paragraph (collapsed,
# same as its original empty/collapsed form).
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$endRange = $gitParaB.Range
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endRange)

Write-Output "done"
